$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header E1 "EQP" -> "Flag3", and add a new header F1 "Flag4"
$ws.Range("E1").Value = "Flag3"
$ws.Range("F1").Value = "Flag4"

# New column F needs blank cells (same formatting as existing row cells) on rows 2-3
$ws.Range("F2").Value = ""
$ws.Range("F3").Value = ""

# Align styles: the previously separate "centered, non-bold" style (used by B2:E3)
# collapses onto the same style as the centered header/body style (style 1),
# so every body cell in row 1 and B2:F3 ends up sharing one consistent centered style.
$ws.Range("B2:F3").HorizontalAlignment = -4108
$ws.Range("F1").HorizontalAlignment = -4108

# Remove the data validation specifically from E1 so the validated header
# range shrinks from A1:E1 to A1:D1 (F1 never had validation to begin with).
$ws.Range("E1").Validation.Delete()

# Update selection to D3
$ws.Range("D3").Select()
